# Add referenced questions when editing a managed form record
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intro")

# Add a description for the "parameters" row heading (D51)
$ws.Range("D51").Value = "Use with a row type of column, space separated parameters from the list below"

# Insert two new blank rows after the existing "parameters" block (before row 54,
# which holds the "action" section heading) so there is room for the new entry.
$ws.Rows("54:55").Insert()

# Document the new "source=question" parameter on row 53
$ws.Range("C53").Value = "source=question"
$ws.Range("D53").Value = "Identify a question that will be the source of data for a column"
